$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 423, shifting the existing rows 423:442 down to 424:443.
$ws.Rows.Item(423).Insert()

# Populate the newly inserted row with the new weekly record.
$ws.Range("A423").Value = 4
$ws.Range("B423").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C423").Value = "Los Lagos"
$ws.Range("D423").Value = 45147
$ws.Range("E423").Value = 10
$ws.Range("F423").Value = 100112032
$ws.Range("G423").Value = "Zapallo italiano"
$ws.Range("H423").Value = "Sin especificar"
$ws.Range("I423").Value = "Primera"
$ws.Range("J423").Value = 70
$ws.Range("K423").Value = 20000
$ws.Range("L423").Value = 20000
$ws.Range("M423").Value = 20000
$ws.Range("N423").Value = "$/caja 50 unidades"
$ws.Range("O423").Value = "Región de Arica y Parinacota"
$ws.Range("P423").Value = 400
$ws.Range("Q423").Value = 50
$ws.Range("R423").Value = "Hortaliza"

Write-Output "applied"
